# Removed Test Case Inter-Dependency
# Update the product name (B1) on both sheets to a unique "-1st" suffixed
# value, and change the short name (B2) on the input sheet from the
# numeric 2565 to the text "256e" so the test case no longer shares
# identifiers with other test runs.

$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "2565-MS-EPP-DB-DL-REC-NOCOM-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-NONE-1st"

# Set B2 (shortname) before B1 (productname) so new shared-string entries
# are appended in the same order as the target workbook.
$wsInput.Range("B2").Value = "256e"
$wsInput.Range("B1").Value = $newProductName

$wsOutput.Range("B1").Value = $newProductName

# Reset the remembered selection on each sheet to B1 (was B15), making
# sure the input sheet ends up as the active/selected tab again.
$wsOutput.Activate()
$wsOutput.Range("B1").Select() | Out-Null
$wsInput.Activate()
$wsInput.Range("B1").Select() | Out-Null
